$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in day-6 (K column) quantities for rows 8-14 ---
# K8 already exists (blank, s=10) - just set its value.
$ws.Range("K8").Value = 74

# K9..K14 don't exist yet; copy the row's existing "day 5" (J) cell
# formatting first so the new cell picks up the same style (s=14)
# that the rest of the row already uses, then overwrite the value.
$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Value = 266

$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("K10").Value = 1050

$ws.Range("J11").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = 298

$ws.Range("J12").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = 57

$ws.Range("J13").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = 72

$ws.Range("J14").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = 61

# --- Update hyperlink display text to be prefixed with "View-source:" ---
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "View-source:" + $hl.Address
}

# T13/T14 carried a redundant duplicate style (identical to T8/T9's);
# re-apply T8's format so the style table collapses the duplicate entry.
$ws.Range("T8").Copy()
$ws.Range("T13").PasteSpecial(-4122)
$ws.Range("T14").PasteSpecial(-4122)

# --- Move the active selection from J14 to K14 ---
[void]$ws.Range("K14").Select()
